# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (market date 2021-11-xx -> serial 44509) for
# "Vega Central Mapocho de Santiago" / Frutilla / Provincia de San Antonio,
# right before the existing row 562. This pushes the existing rows 562:585
# down to 565:588 (dimension grows from A1:T585 to A1:T588) without touching
# their contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 562, shifting everything below down by 3.
$ws.Rows("562:564").Insert()

# Row 562 - Especial
$ws.Range("A562").Value = 9
$ws.Range("B562").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C562").Value = "Metropolitana"
$ws.Range("D562").Value = 44509
$ws.Range("E562").Value = 13
$ws.Range("F562").Value = "Fruta"
$ws.Range("G562").Value = 100101
$ws.Range("H562").Value = "Berries"
$ws.Range("I562").Value = 100112025
$ws.Range("J562").Value = "Frutilla"
$ws.Range("K562").Value = "Sin especificar"
$ws.Range("L562").Value = "Especial"
$ws.Range("M562").Value = 930
$ws.Range("N562").Value = 4500
$ws.Range("O562").Value = 5000
$ws.Range("P562").Value = 4758
$ws.Range("Q562").Value = "`$/bandeja 7 kilos"
$ws.Range("R562").Value = "Provincia de San Antonio"
$ws.Range("S562").Value = 680
$ws.Range("T562").Value = 7

# Row 563 - Primera
$ws.Range("A563").Value = 9
$ws.Range("B563").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C563").Value = "Metropolitana"
$ws.Range("D563").Value = 44509
$ws.Range("E563").Value = 13
$ws.Range("F563").Value = "Fruta"
$ws.Range("G563").Value = 100101
$ws.Range("H563").Value = "Berries"
$ws.Range("I563").Value = 100112025
$ws.Range("J563").Value = "Frutilla"
$ws.Range("K563").Value = "Sin especificar"
$ws.Range("L563").Value = "Primera"
$ws.Range("M563").Value = 660
$ws.Range("N563").Value = 3500
$ws.Range("O563").Value = 4000
$ws.Range("P563").Value = 3689
$ws.Range("Q563").Value = "`$/bandeja 7 kilos"
$ws.Range("R563").Value = "Provincia de San Antonio"
$ws.Range("S563").Value = 527
$ws.Range("T563").Value = 7

# Row 564 - Segunda
$ws.Range("A564").Value = 9
$ws.Range("B564").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C564").Value = "Metropolitana"
$ws.Range("D564").Value = 44509
$ws.Range("E564").Value = 13
$ws.Range("F564").Value = "Fruta"
$ws.Range("G564").Value = 100101
$ws.Range("H564").Value = "Berries"
$ws.Range("I564").Value = 100112025
$ws.Range("J564").Value = "Frutilla"
$ws.Range("K564").Value = "Sin especificar"
$ws.Range("L564").Value = "Segunda"
$ws.Range("M564").Value = 850
$ws.Range("N564").Value = 2500
$ws.Range("O564").Value = 3000
$ws.Range("P564").Value = 2735
$ws.Range("Q564").Value = "`$/bandeja 7 kilos"
$ws.Range("R564").Value = "Provincia de San Antonio"
$ws.Range("S564").Value = 391
$ws.Range("T564").Value = 7

Write-Output "rows inserted"
